$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "KPI Threshold" sheet (sheet1): update a few threshold values
# ---------------------------------------------------------------------------
$wsKpi = $wb.Worksheets.Item("KPI Threshold")
$wsKpi.Range("E14").Value = 0.7
$wsKpi.Range("E15").Value = 98
$wsKpi.Range("E16").Value = 97

# ---------------------------------------------------------------------------
# 2) "Province Area" sheet (sheet4): rebuild the province table, now sorted
#    alphabetically and with several new provinces added, mark "trọng điểm"
#    (priority) provinces with the alternate style, then add an AutoFilter.
# ---------------------------------------------------------------------------
$wsProv = $wb.Worksheets.Item("Province Area")

# Capture style templates into scratch rows below the final table BEFORE we
# start overwriting anything, so every target row can copy from a stable
# source regardless of processing order.
#   - A2:C2 carries cellXf "2"  (non-priority / "khong trong diem" look)
#   - B5:C5 carries cellXf "3"  (priority / "trong diem" look, columns B/C)
$wsProv.Range("A2:C2").Copy($wsProv.Range("A100:C100"))
$wsProv.Range("B5:C5").Copy($wsProv.Range("B101:C101"))

# Full target list (row 1 is the header and is left untouched), sorted
# alphabetically by province code; second element flags "trong diem".
$provinces = @(
    @("BLU", $false),
    @("BPC", $false),
    @("BTN", $false),
    @("CBG", $false),
    @("CMU", $true),
    @("DBN", $false),
    @("DLK", $true),
    @("DNO", $false),
    @("GLI", $false),
    @("HGG", $false),
    @("HNI", $true),
    @("HUG", $false),
    @("KGG", $true),
    @("KTM", $false),
    @("LCI", $false),
    @("LCU", $false),
    @("LDG", $true),
    @("LSN", $false),
    @("NTN", $false),
    @("SLA", $false),
    @("STG", $false),
    @("TNH", $false)
)

$rowNum = 2
foreach ($p in $provinces) {
    $code = $p[0]
    $isPriority = $p[1]

    # Column A always uses the "style 2" look, regardless of priority.
    $wsProv.Range("A100").Copy($wsProv.Range("A$rowNum"))
    $wsProv.Range("A$rowNum").Value = $code

    if ($isPriority) {
        $wsProv.Range("B101:C101").Copy($wsProv.Range("B$rowNum`:C$rowNum"))
    } else {
        $wsProv.Range("B100:C100").Copy($wsProv.Range("B$rowNum`:C$rowNum"))
    }

    $rowNum++
}

# Drop the scratch rows used as copy-source templates.
$wsProv.Range("A100:C101").EntireRow.Delete()

# Apply an AutoFilter over the whole table and register the hidden
# _FilterDatabase defined name the way Excel does internally.
$wsProv.Range("A1:C23").AutoFilter()
$filterName = $wsProv.Names.Add("_xlnm._FilterDatabase", "='Province Area'!`$A`$1:`$C`$23")
$filterName.Visible = $false

# ---------------------------------------------------------------------------
# 3) Restore per-sheet cursor/selection positions to match the edited file.
#    Select in order, finishing on "WorstCell Threshold" so it stays the
#    active tab (matches the unchanged activeTab/tabSelected in the diff).
# ---------------------------------------------------------------------------
$wsKpi.Range("B27").Select()
$wsProv.Range("G22").Select()

$wsWorst = $wb.Worksheets.Item("WorstCell Threshold")
$wsWorst.Range("J26").Select()

Write-Host "edit complete"
